$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(0.9926147277508489, 0.9955398523336033)
    3  = @(0.9936372048519304, 0.9963617723202692)
    4  = @(0.9942998659930995, 0.9968940712668345)
    5  = @(0.9945786998346017, 0.997117960005301)
    6  = @(0.9946255319796338, 0.9971555583673453)
    7  = @(0.9943035907982488, 0.9968970624462087)
    8  = @(0.9929600610674301, 0.995817528259106)
    9  = @(0.9906006454969559, 0.9939188001724441)
    10 = @(0.989033133672735, 0.9926553831429383)
    11 = @(0.988355674866747, 0.9921088820399291)
    12 = @(0.9881042295826724, 0.9919059725120875)
    13 = @(0.9881581567098651, 0.9919494934313052)
    14 = @(0.9883348863814464, 0.9920921077337197)
    15 = @(0.9884438009545853, 0.9921799884222134)
    16 = @(0.9890781214508737, 0.9926916645766087)
    17 = @(0.989476357848556, 0.9930127773699352)
    18 = @(0.9897087662937556, 0.9932001317071769)
    19 = @(0.9897880325774034, 0.9932640239640975)
    20 = @(0.9894336180360679, 0.9929783193494216)
    21 = @(0.9882828385668253, 0.9920501090198105)
    22 = @(0.9875604150241495, 0.9914670000341481)
    23 = @(0.9879432794643023, 0.991776070289318)
    24 = @(0.9894529299347244, 0.9929938892766442)
    25 = @(0.9912096547607049, 0.9944092447426414)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Range("E$row").Value = $pair[0]
    $ws.Range("L$row").Value = $pair[1]
}
